$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new client row (row 6) to the Clients sheet.
# Columns: A=code_client, B=nom, C=contact, D=IFU
# C and D are numeric-looking strings with significant leading zeros, so they
# must be forced to text before assignment to avoid Excel coercing them to
# numbers (which would drop the leading zero).
$ws.Range("A6").Value = "fif012456"
$ws.Range("B6").Value = "Fifi"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "0123456789"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0123456"

# Restore the default "Normal" style on the new numeric-text cells so no
# extra cell formatting is left behind on the new row.
$ws.Range("C6:D6").Style = "Normal"
